$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's worth of entries was recorded (2014-07-29), which needs a new
# row in the data table. Insert a row before the old blank separator row
# (row 164) so every row below (the blank separator + the three summary
# rows) shifts down by one and keeps its relative layout/formulas.
$ws.Rows("164:164").Insert()

# Row 163 (2014-07-28) was left without an end time before; fill it in now
# and let the time-spent formulas compute from it.
$ws.Range("E163").Value = 0.73958333333333337
$ws.Range("F163").Formula = "=(E163-D163)*24*60"
$ws.Range("G163").Formula = "=F163/60"

# Row 164 (2014-07-29): newly recorded day.
$ws.Range("A164").Value = 2014
$ws.Range("B164").Value = 7
$ws.Range("C164").Value = 29
$ws.Range("D164").Value = 0.33333333333333331
$ws.Range("E164").Value = 0.5
$ws.Range("F164").Formula = "=(E164-D164)*24*60"
$ws.Range("G164").Formula = "=F164/60"

# The summary formulas (now on rows 166-168 after the insert) need their
# ranges extended to cover the new data row.
$ws.Range("F166").Formula = "=SUM(F2:F164)"
$ws.Range("F167").Formula = "=F166/60"
$ws.Range("F168").Formula = "=F167/38.5"

$ws.Range("F164").Select() | Out-Null
